# Update countries & provincias Spain
#
# The source COVID dashboard re-ranked a handful of countries (their total
# case counts overtook neighbours in the list) and refreshed the day's
# figures. Three rows swap the country they display (the newcomer gets
# fresh numbers, the displaced country keeps its old numbers but drops one
# row), a batch of other rows simply get updated totals, and the "data as
# of" timestamp moves from 08:40 to 09:57.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 09:57"

# Rusia (row 6) - refreshed totals
$ws.Range("B6").Value = 606881
$ws.Range("C6").Value = 7176
$ws.Range("D6").Value = 368822
$ws.Range("E6").Value = 229546
$ws.Range("G6").Value = 154
$ws.Range("H6").Value = 8513

# Singapur (row 35) - refreshed totals
$ws.Range("B35").Value = 42623
$ws.Range("C35").Value = 191
$ws.Range("E35").Value = 6602

# Afganistan (row 44) - refreshed totals
$ws.Range("B44").Value = 29640
$ws.Range("C44").Value = 159
$ws.Range("D44").Value = 9869
$ws.Range("E44").Value = 19132
$ws.Range("G44").Value = 21
$ws.Range("H44").Value = 639

# Armenia overtakes Israel and Nigeria: row 51 becomes Armenia (new data),
# row 52 becomes Israel (its former row-51 data), row 53 becomes Nigeria
# (its former row-52 data).
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 21717
$ws.Range("C51").Value = 711
$ws.Range("D51").Value = 10797
$ws.Range("E51").Value = 10534
$ws.Range("G51").Value = 14
$ws.Range("H51").Value = 386

$ws.Range("A52").Value = "Israel"
$ws.Range("B52").Value = 21512
$ws.Range("D52").Value = 15869
$ws.Range("E52").Value = 5335
$ws.Range("H52").Value = 308

$ws.Range("A53").Value = "Nigeria"
$ws.Range("B53").Value = 21371
$ws.Range("D53").Value = 7338
$ws.Range("E53").Value = 13500
$ws.Range("H53").Value = 533

# Chequia (row 67) - refreshed totals
$ws.Range("B67").Value = 10651
$ws.Range("C67").Value = 1
$ws.Range("D67").Value = 7559
$ws.Range("E67").Value = 2752
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 340

# El Salvador overtakes Guinea: row 82 becomes El Salvador (new data),
# row 83 becomes Guinea (its former row-82 data).
$ws.Range("A82").Value = "El Salvador"
$ws.Range("B82").Value = 5150
$ws.Range("C82").Value = 177
$ws.Range("D82").Value = 2924
$ws.Range("E82").Value = 2107
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 119

$ws.Range("A83").Value = "Guinea"
$ws.Range("B83").Value = 5040
$ws.Range("D83").Value = 3685
$ws.Range("E83").Value = 1327
$ws.Range("H83").Value = 28

# Hungria (row 91) - refreshed totals
$ws.Range("B91").Value = 4114
$ws.Range("C91").Value = 7
$ws.Range("D91").Value = 2618
$ws.Range("E91").Value = 920
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 576

# Somalia (row 98) - refreshed totals
$ws.Range("B98").Value = 2835
$ws.Range("C98").Value = 23
$ws.Range("D98").Value = 829
$ws.Range("E98").Value = 1916

# Estonia (row 108) - refreshed totals
$ws.Range("B108").Value = 1983
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 1783
$ws.Range("E108").Value = 131

# Eslovaquia (row 116) - refreshed totals
$ws.Range("B116").Value = 1607
$ws.Range("C116").Value = 18
$ws.Range("E116").Value = 131

# Georgia (row 133) - refreshed totals
$ws.Range("D133").Value = 771
$ws.Range("E133").Value = 129

# Montenegro (row 156) - refreshed totals
$ws.Range("B156").Value = 383
$ws.Range("C156").Value = 5
$ws.Range("E156").Value = 59

# Seychelles overtakes Montserrat: row 211 becomes Seychelles, row 212
# becomes Montserrat (values simply swap between the two rows).
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = "Montserrat"
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
